$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44208
$ws.Cells.Item(2, 11).Value = 'Lapins'
$ws.Cells.Item(2, 12).Value = 'Segunda'
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 10500
$ws.Cells.Item(2, 15).Value = 11000
$ws.Cells.Item(2, 16).Value = 10750
$ws.Cells.Item(2, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(2, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(2, 19).Value = 896
$ws.Cells.Item(2, 20).Value = 12

$ws.Cells.Item(3, 4).Value = 44901
$ws.Cells.Item(3, 11).Value = 'Bing'
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 500
$ws.Cells.Item(3, 14).Value = 12000
$ws.Cells.Item(3, 15).Value = 13000
$ws.Cells.Item(3, 16).Value = 12500
$ws.Cells.Item(3, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(3, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 19).Value = 833
$ws.Cells.Item(3, 20).Value = 15

$ws.Cells.Item(4, 4).Value = 44901
$ws.Cells.Item(4, 11).Value = 'Lapins'
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 500
$ws.Cells.Item(4, 14).Value = 12000
$ws.Cells.Item(4, 15).Value = 13000
$ws.Cells.Item(4, 16).Value = 12500
$ws.Cells.Item(4, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(4, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(4, 19).Value = 833
$ws.Cells.Item(4, 20).Value = 15

$ws.Cells.Item(5, 4).Value = 44917
$ws.Cells.Item(5, 11).Value = 'Bing'
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 400
$ws.Cells.Item(5, 14).Value = 5000
$ws.Cells.Item(5, 15).Value = 6000
$ws.Cells.Item(5, 16).Value = 5625
$ws.Cells.Item(5, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(5, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(5, 19).Value = 562
$ws.Cells.Item(5, 20).Value = 10

$ws.Cells.Item(6, 4).Value = 44917
$ws.Cells.Item(6, 11).Value = 'Santina'
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 400
$ws.Cells.Item(6, 14).Value = 5000
$ws.Cells.Item(6, 15).Value = 6000
$ws.Cells.Item(6, 16).Value = 5500
$ws.Cells.Item(6, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(6, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(6, 19).Value = 550
$ws.Cells.Item(6, 20).Value = 10

$ws.Cells.Item(7, 4).Value = 44580
$ws.Cells.Item(7, 11).Value = 'Sweet Heart'
$ws.Cells.Item(7, 12).Value = 'Segunda'
$ws.Cells.Item(7, 13).Value = 300
$ws.Cells.Item(7, 14).Value = 7000
$ws.Cells.Item(7, 15).Value = 8000
$ws.Cells.Item(7, 16).Value = 7500
$ws.Cells.Item(7, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(7, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(7, 19).Value = 750
$ws.Cells.Item(7, 20).Value = 10

$ws.Cells.Item(8, 4).Value = 44594
$ws.Cells.Item(8, 11).Value = 'Santina'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 160
$ws.Cells.Item(8, 14).Value = 5000
$ws.Cells.Item(8, 15).Value = 6000
$ws.Cells.Item(8, 16).Value = 5500
$ws.Cells.Item(8, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 1100
$ws.Cells.Item(8, 20).Value = 5

$ws.Cells.Item(9, 4).Value = 44537
$ws.Cells.Item(9, 11).Value = 'Brooks'
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 29000
$ws.Cells.Item(9, 15).Value = 30000
$ws.Cells.Item(9, 16).Value = 29500
$ws.Cells.Item(9, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(9, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(9, 19).Value = 1475
$ws.Cells.Item(9, 20).Value = 20

$ws.Cells.Item(10, 4).Value = 44908
$ws.Cells.Item(10, 11).Value = 'Rainier'
$ws.Cells.Item(10, 12).Value = 'Segunda'
$ws.Cells.Item(10, 13).Value = 250
$ws.Cells.Item(10, 14).Value = 15000
$ws.Cells.Item(10, 15).Value = 16000
$ws.Cells.Item(10, 16).Value = 15600
$ws.Cells.Item(10, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10, 19).Value = 1560
$ws.Cells.Item(10, 20).Value = 10

$ws.Cells.Item(11, 4).Value = 44931
$ws.Cells.Item(11, 11).Value = 'Lapins'
$ws.Cells.Item(11, 12).Value = 'Segunda'
$ws.Cells.Item(11, 13).Value = 250
$ws.Cells.Item(11, 14).Value = 6000
$ws.Cells.Item(11, 15).Value = 6500
$ws.Cells.Item(11, 16).Value = 6250
$ws.Cells.Item(11, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(11, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(11, 19).Value = 625
$ws.Cells.Item(11, 20).Value = 10

$ws.Cells.Item(12, 4).Value = 44931
$ws.Cells.Item(12, 11).Value = 'Lapins'
$ws.Cells.Item(12, 12).Value = 'Segunda'
$ws.Cells.Item(12, 13).Value = 400
$ws.Cells.Item(12, 14).Value = 3000
$ws.Cells.Item(12, 15).Value = 3300
$ws.Cells.Item(12, 16).Value = 3150
$ws.Cells.Item(12, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(12, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 19).Value = 630
$ws.Cells.Item(12, 20).Value = 5

$ws.Cells.Item(13, 4).Value = 44161
$ws.Cells.Item(13, 11).Value = 'Bing'
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 160
$ws.Cells.Item(13, 14).Value = 39000
$ws.Cells.Item(13, 15).Value = 40000
$ws.Cells.Item(13, 16).Value = 39500
$ws.Cells.Item(13, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(13, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(13, 19).Value = 1975
$ws.Cells.Item(13, 20).Value = 20

$ws.Cells.Item(14, 4).Value = 44571
$ws.Cells.Item(14, 11).Value = 'Brooks'
$ws.Cells.Item(14, 12).Value = 'Segunda'
$ws.Cells.Item(14, 13).Value = 400
$ws.Cells.Item(14, 14).Value = 8500
$ws.Cells.Item(14, 15).Value = 9000
$ws.Cells.Item(14, 16).Value = 8750
$ws.Cells.Item(14, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(14, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(14, 19).Value = 875
$ws.Cells.Item(14, 20).Value = 10

$ws.Cells.Item(15, 4).Value = 44532
$ws.Cells.Item(15, 11).Value = 'Brooks'
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 400
$ws.Cells.Item(15, 14).Value = 27000
$ws.Cells.Item(15, 15).Value = 28000
$ws.Cells.Item(15, 16).Value = 27500
$ws.Cells.Item(15, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(15, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(15, 19).Value = 2292
$ws.Cells.Item(15, 20).Value = 12

$ws.Cells.Item(16, 4).Value = 44557
$ws.Cells.Item(16, 11).Value = 'Lapins'
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 250
$ws.Cells.Item(16, 14).Value = 9000
$ws.Cells.Item(16, 15).Value = 10000
$ws.Cells.Item(16, 16).Value = 9500
$ws.Cells.Item(16, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(16, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(16, 19).Value = 950
$ws.Cells.Item(16, 20).Value = 10

$ws.Cells.Item(17, 4).Value = 44568
$ws.Cells.Item(17, 11).Value = 'Santina'
$ws.Cells.Item(17, 12).Value = 'Segunda'
$ws.Cells.Item(17, 13).Value = 200
$ws.Cells.Item(17, 14).Value = 15000
$ws.Cells.Item(17, 15).Value = 16000
$ws.Cells.Item(17, 16).Value = 15500
$ws.Cells.Item(17, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(17, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(17, 19).Value = 1292
$ws.Cells.Item(17, 20).Value = 12

$ws.Cells.Item(18, 4).Value = 44922
$ws.Cells.Item(18, 11).Value = 'Bing'
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 13).Value = 200
$ws.Cells.Item(18, 14).Value = 5000
$ws.Cells.Item(18, 15).Value = 6000
$ws.Cells.Item(18, 16).Value = 5500
$ws.Cells.Item(18, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(18, 18).Value = 'Región del Maule'
$ws.Cells.Item(18, 19).Value = 550
$ws.Cells.Item(18, 20).Value = 10

$ws.Cells.Item(19, 4).Value = 44229
$ws.Cells.Item(19, 11).Value = 'Santina'
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 250
$ws.Cells.Item(19, 14).Value = 6500
$ws.Cells.Item(19, 15).Value = 7000
$ws.Cells.Item(19, 16).Value = 6750
$ws.Cells.Item(19, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(19, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(19, 19).Value = 1350
$ws.Cells.Item(19, 20).Value = 5

$ws.Cells.Item(20, 4).Value = 44210
$ws.Cells.Item(20, 11).Value = 'Rainier'
$ws.Cells.Item(20, 12).Value = 'Segunda'
$ws.Cells.Item(20, 13).Value = 250
$ws.Cells.Item(20, 14).Value = 21000
$ws.Cells.Item(20, 15).Value = 22000
$ws.Cells.Item(20, 16).Value = 21500
$ws.Cells.Item(20, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(20, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(20, 19).Value = 1194
$ws.Cells.Item(20, 20).Value = 18

$ws.Cells.Item(21, 4).Value = 44943
$ws.Cells.Item(21, 11).Value = 'Santina'
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 600
$ws.Cells.Item(21, 14).Value = 14000
$ws.Cells.Item(21, 15).Value = 15000
$ws.Cells.Item(21, 16).Value = 14333
$ws.Cells.Item(21, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(21, 18).Value = 'Región del Maule'
$ws.Cells.Item(21, 19).Value = 956
$ws.Cells.Item(21, 20).Value = 15

$ws.Cells.Item(22, 4).Value = 44921
$ws.Cells.Item(22, 11).Value = 'Bing'
$ws.Cells.Item(22, 12).Value = 'Primera'
$ws.Cells.Item(22, 13).Value = 320
$ws.Cells.Item(22, 14).Value = 7500
$ws.Cells.Item(22, 15).Value = 8000
$ws.Cells.Item(22, 16).Value = 7781
$ws.Cells.Item(22, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(22, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(22, 19).Value = 778
$ws.Cells.Item(22, 20).Value = 10

$ws.Cells.Item(23, 4).Value = 44914
$ws.Cells.Item(23, 11).Value = 'Brooks'
$ws.Cells.Item(23, 12).Value = 'Primera'
$ws.Cells.Item(23, 13).Value = 700
$ws.Cells.Item(23, 14).Value = 7000
$ws.Cells.Item(23, 15).Value = 8000
$ws.Cells.Item(23, 16).Value = 7429
$ws.Cells.Item(23, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(23, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(23, 19).Value = 743
$ws.Cells.Item(23, 20).Value = 10

$ws.Cells.Item(24, 4).Value = 44914
$ws.Cells.Item(24, 11).Value = 'Lapins'
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 550
$ws.Cells.Item(24, 14).Value = 7000
$ws.Cells.Item(24, 15).Value = 8000
$ws.Cells.Item(24, 16).Value = 7455
$ws.Cells.Item(24, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(24, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(24, 19).Value = 746
$ws.Cells.Item(24, 20).Value = 10

$ws.Cells.Item(25, 4).Value = 44175
$ws.Cells.Item(25, 11).Value = 'Rainier'
$ws.Cells.Item(25, 12).Value = 'Segunda'
$ws.Cells.Item(25, 13).Value = 270
$ws.Cells.Item(25, 14).Value = 25000
$ws.Cells.Item(25, 15).Value = 26000
$ws.Cells.Item(25, 16).Value = 25500
$ws.Cells.Item(25, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(25, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(25, 19).Value = 1417
$ws.Cells.Item(25, 20).Value = 18
